$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark from its old location -------
# (Word keeps only one "_GoBack" bookmark at a time, marking the last edit
#  location; it is about to be re-created at the new edit position below.)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Locate the sentence that needs the "." inserted --------------------
$target = $d.Content
$target.Find.Execute("Reached 103852 scores in T-Rex Game on Chrome browser.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $target.Find.Found) {
    throw "Could not locate the T-Rex high score sentence."
}

$sentenceStart = $target.Start
$splitPos = $sentenceStart + 11   # right after "Reached 103", before "852..."

# Guard the left edge of the run (boundary with the preceding "  " run) so
# that the upcoming edit does not coalesce it with the "Reached 103..." run.
$leftGuard = $d.Range($splitPos - 11, $splitPos - 11)
$d.Bookmarks.Add("ZZ_GUARD_LEFT", $leftGuard)

# --- 3. Type the "." in the middle of the run -------------------------------
$insertionPoint = $d.Range($splitPos, $splitPos)
$insertionPoint.InsertBefore(".")

# Nudge the freshly-typed "." with a no-op formatting round trip so it is
# kept as its own run instead of being re-merged with "Reached 103".
$dotRange = $d.Range($splitPos, $splitPos + 1)
$dotRange.Font.Bold = 1
$dotRange.Font.Bold = 0

# --- 4. Drop the "_GoBack" bookmark exactly where the user stopped typing --
$newGoBack = $d.Range($splitPos + 1, $splitPos + 1)
$d.Bookmarks.Add("_GoBack", $newGoBack)

# --- 5. Remove the temporary guard bookmark --------------------------------
$d.Bookmarks("ZZ_GUARD_LEFT").Delete()
